$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1140
$ws.Range("F6").Value = 3393
$ws.Range("F10").Value = 597
$ws.Range("F12").Value = 155
$ws.Range("F14").Value = 1816
$ws.Range("F16").Value = 392
$ws.Range("F19").Value = 682
$ws.Range("F22").Value = 794
$ws.Range("F23").Value = 79977
$ws.Range("F24").Value = 79977
$ws.Range("F27").Value = 33833
$ws.Range("F28").Value = 33834
$ws.Range("F33").Value = 50
$ws.Range("F38").Value = 2517
$ws.Range("F39").Value = 2517
$ws.Range("F40").Value = 1211
$ws.Range("F41").Value = 5499
$ws.Range("F42").Value = 792
$ws.Range("F47").Value = 414

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 1976
$ws.Range("F42").Value = 35
$ws.Range("F47").Value = 194

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 148

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 3393
$ws.Range("F14").Value = 597
$ws.Range("F17").Value = 148
$ws.Range("F18").Value = 1816
$ws.Range("F25").Value = 682
$ws.Range("F29").Value = 79977
$ws.Range("F31").Value = 33834
$ws.Range("F36").Value = 50
$ws.Range("F45").Value = 2517
$ws.Range("F46").Value = 1211
$ws.Range("F47").Value = 792
$ws.Range("F51").Value = 35
$ws.Range("F55").Value = 194
